$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pepID")

# Find the last used row in column A and append the new PatientID right after it
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "PEP_ID-2009435"
